# Update Name of Algo
# Applies the updated KNN-imputed values to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B4").Value = 7.239
$ws.Range("D5").Value = -8.251999999999999
$ws.Range("A8").Value = -21.54799999999999
$ws.Range("D8").Value = -7.922
$ws.Range("A10").Value = -21.794
$ws.Range("B11").Value = 6.578
$ws.Range("A12").Value = -21.178
$ws.Range("B12").Value = 6.245
$ws.Range("D12").Value = -6.483
$ws.Range("D13").Value = -8.111000000000001
$ws.Range("B15").Value = 6.17
$ws.Range("D15").Value = -8.193
$ws.Range("B17").Value = 4.851000000000001
$ws.Range("A18").Value = -21.739
$ws.Range("D21").Value = -8.395999999999999
$ws.Range("A25").Value = -21.688
$ws.Range("D25").Value = -7.887
$ws.Range("B26").Value = 5.999000000000001
$ws.Range("B27").Value = 5.994999999999999
$ws.Range("B28").Value = 6.257999999999999
$ws.Range("B32").Value = 5.839
$ws.Range("D32").Value = -7.756
$ws.Range("D36").Value = -7.640000000000001
$ws.Range("A37").Value = -20.964
$ws.Range("B37").Value = 8.109
$ws.Range("D38").Value = -7.752
$ws.Range("B41").Value = 8.386999999999999
$ws.Range("D41").Value = -7.927
$ws.Range("B47").Value = 5.69
$ws.Range("D50").Value = -8.151
$ws.Range("B51").Value = 6.991
$ws.Range("D52").Value = -7.528
$ws.Range("A55").Value = -22.075
$ws.Range("D59").Value = -8.109
$ws.Range("B65").Value = 6.511
$ws.Range("D67").Value = -7.556
$ws.Range("A68").Value = -21.483
$ws.Range("B73").Value = 6.022
$ws.Range("A77").Value = -21.036
$ws.Range("A78").Value = -20.572
$ws.Range("A79").Value = -21.017
$ws.Range("A80").Value = -21.215
$ws.Range("A81").Value = -21.738
$ws.Range("A82").Value = -21.563
$ws.Range("A84").Value = -21.218
$ws.Range("B84").Value = 7.855
$ws.Range("D84").Value = -8.145999999999999
$ws.Range("B85").Value = 6.017
$ws.Range("D86").Value = -8.199999999999999
$ws.Range("D88").Value = -8.242999999999999
$ws.Range("B89").Value = 4.705
$ws.Range("D89").Value = -7.798999999999999
$ws.Range("B93").Value = 6.077
$ws.Range("B95").Value = 5.781000000000001
$ws.Range("D95").Value = -7.918000000000001
$ws.Range("B98").Value = 6.752
$ws.Range("B99").Value = 6.06
$ws.Range("A101").Value = -21.946
$ws.Range("B101").Value = 5.936
$ws.Range("A102").Value = -21.007
$ws.Range("B102").Value = 7.039
$ws.Range("D105").Value = -7.961999999999999
